# Add a "Save" column (H) to the s_vals sheet, matching the header style
# used by the existing columns (B1:G1) and appending a numeric 0 value
# for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same shared cell style (bold font,
# border, centered/top alignment) instead of creating a near-duplicate style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for the single data row, numeric 0.
$ws.Range("H2").Value = 0
